{"js": "const replacements = [\n  [\"89\u00d766=5874\", \"60\u00d795=5700\"],\n  [\"71\u00d780=5680\", \"70\u00d722=1540\"],\n  [\"82\u00d733=2706\", \"42\u00d790=3780\"],\n  [\"23\u00d776=1748\", \"98\u00d789=8722\"],\n  [\"41\u00d745=1845\", \"30\u00d749=1470\"],\n  [\"70\u00d741=2870\", \"85\u00d782=6970\"],\n  [\"15\u00d711=165\", \"49\u00d722=1078\"],\n  [\"98\u00d752=5096\", \"11\u00d782=902\"],\n  [\"33\u00d765=2145\", \"44\u00d725=1100\"],\n  [\"73\u00d726=1898\", \"57\u00d779=4503\"],\n  [\"22\u00d712=264\", \"35\u00d754=1890\"],\n  [\"61\u00d737=2257\", \"12\u00d727=324\"],\n  [\"13\u00d766=858\", \"93\u00d759=5487\"],\n  [\"76\u00d790=6840\", \"14\u00d754=756\"],\n  [\"65\u00d768=4420\", \"56\u00d719=1064\"],\n  [\"97\u00d759=5723\", \"34\u00d730=1020\"],\n  [\"89\u00d773=6497\", \"56\u00d778=4368\"],\n  [\"70\u00d792=6440\", \"35\u00d724=840\"],\n  [\"42\u00d772=3024\", \"13\u00d792=1196\"],\n  [\"16\u00d795=1520\", \"55\u00d726=1430\"],\n  [\"97\u00d750=4850\", \"24\u00d720=480\"],\n  [\"87\u00d767=5829\", \"44\u00d792=4048\"],\n  [\"74\u00d752=3848\", \"39\u00d784=3276\"],\n  [\"50\u00d753=2650\", \"75\u00d739=2925\"],\n  [\"53\u00d788=4664\", \"82\u00d792=7544\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"89\u00d766=5874\", \"60\u00d795=5700\"),\n    @(\"71\u00d780=5680\", \"70\u00d722=1540\"),\n    @(\"82\u00d733=2706\", \"42\u00d790=3780\"),\n    @(\"23\u00d776=1748\", \"98\u00d789=8722\"),\n    @(\"41\u00d745=1845\", \"30\u00d749=1470\"),\n    @(\"70\u00d741=2870\", \"85\u00d782=6970\"),\n    @(\"15\u00d711=165\", \"49\u00d722=1078\"),\n    @(\"98\u00d752=5096\", \"11\u00d782=902\"),\n    @(\"33\u00d765=2145\", \"44\u00d725=1100\"),\n    @(\"73\u00d726=1898\", \"57\u00d779=4503\"),\n    @(\"22\u00d712=264\", \"35\u00d754=1890\"),\n    @(\"61\u00d737=2257\", \"12\u00d727=324\"),\n    @(\"13\u00d766=858\", \"93\u00d759=5487\"),\n    @(\"76\u00d790=6840\", \"14\u00d754=756\"),\n    @(\"65\u00d768=4420\", \"56\u00d719=1064\"),\n    @(\"97\u00d759=5723\", \"34\u00d730=1020\"),\n    @(\"89\u00d773=6497\", \"56\u00d778=4368\"),\n    @(\"70\u00d792=6440\", \"35\u00d724=840\"),\n    @(\"42\u00d772=3024\", \"13\u00d792=1196\"),\n    @(\"16\u00d795=1520\", \"55\u00d726=1430\"),\n    @(\"97\u00d750=4850\", \"24\u00d720=480\"),\n    @(\"87\u00d767=5829\", \"44\u00d792=4048\"),\n    @(\"74\u00d752=3848\", \"39\u00d784=3276\"),\n    @(\"50\u00d753=2650\", \"75\u00d739=2925\"),\n    @(\"53\u00d788=4664\", \"82\u00d792=7544\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
